$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newOrders = @(
    @{ Row = 7;  OrderId = 6; Items = "[1, 1, 2, 2]" },
    @{ Row = 8;  OrderId = 7; Items = "[2, 2, 2]" },
    @{ Row = 9;  OrderId = 8; Items = "[7, 7, 5]" },
    @{ Row = 10; OrderId = 9; Items = "[1, 1, 4]" }
)

foreach ($order in $newOrders) {
    $r = $order.Row
    $ws.Cells.Item($r, 1).Value = $order.OrderId
    $ws.Cells.Item($r, 2).Value = "dineIn"
    $ws.Cells.Item($r, 3).Value = $order.Items
    $ws.Cells.Item($r, 4).Value = $false
    $ws.Cells.Item($r, 5).Value = "InProgress"
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 7).Value = 7
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 9).Value = 0
}
